$d = $word.ActiveDocument

# 1. Title heading: shorten title (also replaces the second occurrence near the bottom)
$d.Content.Find.Execute(
    "Play Chilli Pop for Free: Exciting Mexican-Themed Slot Game",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Chilli Pop for Free", 2) | Out-Null

# 2. "What we like" bullet: graphics and visual design
$d.Content.Find.Execute(
    "High-quality graphics and visual design",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "High-quality graphics and visually stunning design", 2) | Out-Null

# 3. "What we don't like" bullet: volatility -> payout frequency
$d.Content.Find.Execute(
    "Low to medium volatility",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Frequency of payouts could be higher", 2) | Out-Null

# 4. "What we don't like" bullet: minimum bet wording
$d.Content.Find.Execute(
    "Minimum bet of €0.50 per spin",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Minimum bet amount is €0.50 per spin", 2) | Out-Null

# 5. Meta description (italic run) replacement
$d.Content.Find.Execute(
    "Read our review of Chilli Pop, a colorful Mexican cuisine-inspired online slot game. Play for free and enjoy free spins, bonus features, and multipliers.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Experience the Mexican-themed slot game Chilli Pop for free and enjoy cluster-based winning combinations.", 2) | Out-Null
